$wb = $excel.ActiveWorkbook

# Rename "Paineis DARQ" -> "PAINEIS DARQ"
$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"

# Rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"

# Delete "Desarquivamentos Pendentes" sheet
$excel.DisplayAlerts = $false
[void]$wb.Worksheets.Item("Desarquivamentos Pendentes").Delete()
$excel.DisplayAlerts = $true

# Deleting the last (and previously-selected) sheet shifts the active tab;
# restore "PAINEIS DARQ" (originally the active/selected sheet) as active.
[void]$wb.Worksheets.Item("PAINEIS DARQ").Activate()
